$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.141.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.54%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.814.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.70%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.68%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'232.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.19%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.78%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'41.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.91%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +6.43%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0685"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.15%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0999"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.25%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.079.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.60%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.817.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.57%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'11.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.94%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.03%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'4.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.27%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'35.098.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.70%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'69.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.30%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.25%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'239.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.22%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.54%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.63%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.08%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'172.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.42%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.54%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'17.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.14%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.22%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +19.36%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.73%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.332.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'4.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.81%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +3.71%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.77%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -7.92%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +4.30%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'92.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.04%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.683"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.74%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.24%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Maker"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'1.310.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.57%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'ARBITRUM"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.51%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'WEMIXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'1.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.48%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.71%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +1.49%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -5.18%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.53%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'6.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.79%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.29%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.992.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.96%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.63%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0647"
$ws.Range("D51").Style = "Normal"
